$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3.5
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.9
$ws.Range("W2").Value = 1.77
$ws.Range("X2").Value = 1.87
$ws.Range("W3").Value = 1.8
$ws.Range("X3").Value = 1.8
$ws.Range("G4").Value = 1.6
$ws.Range("H4").Value = 3.8
$ws.Range("I4").Value = 4.75
$ws.Range("J4").Value = 2.2
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 4
$ws.Range("Q4").Value = 1.73
$ws.Range("R4").Value = 2.08
$ws.Range("S4").Value = 2.75
$ws.Range("T4").Value = 1.4
$ws.Range("Z4").Value = 8.5
$ws.Range("AC4").Value = 13
$ws.Range("AI4").Value = 15
$ws.Range("AJ4").Value = 26
$ws.Range("AK4").Value = 15
$ws.Range("G5").Value = 1.04
$ws.Range("H5").Value = 8.25
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 1.26
$ws.Range("K5").Value = 3.45
$ws.Range("L5").Value = 22
$ws.Range("Q5").Value = 1.24
$ws.Range("R5").Value = 3.8
$ws.Range("S5").Value = 1.62
$ws.Range("T5").Value = 2.02
$ws.Range("Y5").Value = 9.25
$ws.Range("Z5").Value = 5.6
$ws.Range("AA5").Value = 14
$ws.Range("AB5").Value = 4.8
$ws.Range("AC5").Value = 11.75
$ws.Range("AD5").Value = 50
$ws.Range("AF5").Value = 24
$ws.Range("AG5").Value = 60
$ws.Range("AH5").Value = 300
$ws.Range("AI5").Value = 110
$ws.Range("AK5").Value = 150
$ws.Range("AN5").Value = 500
$ws.Range("G6").Value = 1.91
$ws.Range("I6").Value = 4.1
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 12
$ws.Range("O6").Value = 1.22
$ws.Range("T6").Value = 1.4
$ws.Range("AB6").Value = 17
$ws.Range("AJ6").Value = 21
$ws.Range("AK6").Value = 13
$ws.Range("O7").Value = 1.06
$ws.Range("P7").Value = 10
$ws.Range("G8").Value = 2.75
$ws.Range("I8").Value = 2.45
$ws.Range("J8").Value = 3.25
$ws.Range("L8").Value = 3
$ws.Range("AA8").Value = 11
$ws.Range("AB8").Value = 29
$ws.Range("AL8").Value = 23
$ws.Range("G9").Value = 1.29
$ws.Range("H9").Value = 6
$ws.Range("I9").Value = 8
$ws.Range("J9").Value = 1.73
$ws.Range("L9").Value = 7.5
$ws.Range("Q9").Value = 1.4
$ws.Range("R9").Value = 2.88
$ws.Range("U9").Value = 1.22
$ws.Range("V9").Value = 4
$ws.Range("W9").Value = 1.8
$ws.Range("X9").Value = 1.91
$ws.Range("Y9").Value = 10
$ws.Range("AB9").Value = 9
$ws.Range("AO9").Value = 201
